# Annual Vehicle Maint Cost.xlsx - transportation sector calibration
# Redesignate "passenger ships" row on the AVMC-passenger sheet so that it
# mirrors the LDVs ("taxis") row instead of using the Cost Data ships figures.

$wb = $excel.ActiveWorkbook

$wsPass = $wb.Worksheets.Item("AVMC-passenger")

# Row 6 = "ships". Point every column at the corresponding cell in row 2
# ("LDVs") instead of the old Cost Data-driven ships formulas / hard zeros.
$wsPass.Range("B6:H6").ClearFormats()
$wsPass.Range("B6:H6").NumberFormat = "0"

$wsPass.Range("B6").Formula = "=B2"
$wsPass.Range("C6:H6").Formula = "=C2"

# --- Cosmetic cursor/selection bookkeeping to mirror the saved view state ---
$wsPass.Activate()
$wsPass.Range("C10").Select()

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
$wsAbout.Range("B70").Select()
